$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Note 2 text in A19 with the new, expanded note.
$ws.Range("A19").Value = "Note 2: The first 17 largest thermal generators account for 92.46% of all thermal generation capacity. 4 generators ( CG1, CG2, CG3, and CG4) with capacity 181MW are not dispatched because of abnormally high prices 15892`$/MWh. Therefore, they are not modeled as agents."

# Update the active cell selection to match the new state.
$ws.Range("B23").Select()
